$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column C for every existing data
#    row (2-85) from 2023-09-23 (45192) to 2023-10-03 (45202).
for ($r = 2; $r -le 85; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2. Row 85 picks up an explicit standard row height (15pt / customHeight)
#    in the new file, matching the rest of the data rows.
$ws.Rows.Item(85).RowHeight = 15

# 3. Append the new record as row 86.
$ws.Range("A86").Value = "A 46447-2023"
$ws.Range("B86").Value = 45197
$ws.Range("C86").Value = 45202
$ws.Range("B86:C86").NumberFormat = "YYYY-MM-DD"
$ws.Range("D86").Value = "VÄRMLANDS LÄN"
$ws.Range("E86").Value = "MUNKFORS"
$ws.Range("F86").Value = "Bergvik skog väst AB"
$ws.Range("G86").Value = 1.4
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("P86").Value = 0
$ws.Range("Q86").Value = 0
$ws.Range("R86").Value = ""
$ws.Range("R86").WrapText = $true
